# Applies the cryptos-list update described by the commit diff.
# Each target cell is explicitly formatted as Text ("@") before the
# literal value is written, so numeric-looking strings (e.g. "19.50",
# "1.00", "208.48") keep their exact text form instead of being
# coerced into floating point numbers (which would drop trailing
# zeros / renormalize the value) -- matching the original inline-string
# cells in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.110.18'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.571.72'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.48'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.59%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -3.68%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.50'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.795.30'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.79%  '
$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.04'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.559.32'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.511'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.23'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.117.80'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0725'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.25'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.55%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '206.56'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.83%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.84'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.13'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.95'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.85%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -2.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.19'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.79%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.24%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.278.30'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.610'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.53%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.49%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -7.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.812'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.50%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.66%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.48%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.761'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.24'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.708.54'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.06'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₆0105'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.51'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.54%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.28%  '
